# ITSADSSD-21811 - Changes to add quotes in output file and error handling
#
# This script reproduces, via the Excel COM object model, the edits made to
# Processes/PBI_LocationCheck/Data/Config.xlsx:
#   - Settings sheet: new "Output_Header" row (CSV header for output file),
#     NumberOfEmails value bumped from 10 to 30.
#   - Constants sheet: removed the email-sender name/address rows, added a
#     "NoEmail_Message" row, trimmed the response-email body text, split the
#     old "Email_Alert_MessageBody" row into distinct "EmailMsg_InitError"
#     and "EmailMsg_MoveFolder" rows, and populated "Email_Alert_Subject".
#   - Assets sheet: removed the UQ_SMTP_SERVER / UQ_SMTP_PORT asset rows.
#   - Constants becomes the active sheet/tab; selections on each sheet are
#     updated to reflect where the editor was last working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Remove any pre-existing hyperlinks so they can be re-added at their
# (possibly shifted) final locations further down.
$settings.Range("A1").Hyperlinks.Delete()

# Insert a new row above the old "logF_BusinessProcessName" row (row 5) to
# hold the new Output_Header setting; everything below shifts down by one.
$settings.Rows.Item(5).Insert()

$settings.Range("A5").Value = "Output_Header"
$settings.Range("B5").Value = "EMPLID,FAMILY_NAME,GIVEN_NAMES,BIRTH_DATE,PASSPORT_NUMBER,PASSPORT_COUNTRY,LOCATION,UPDATETIMESTAMP,ERROR"
$settings.Range("C5").Value = "Header fields for output file "

# NumberOfEmails value changes from 10 to 30 (now on row 9 after the insert).
$settings.Range("B9").Value = 30

# Re-create the two hyperlinks at their shifted cell locations.
$settings.Hyperlinks.Add($settings.Range("B8"), "https://outlook.office365.com/EWS/Exchange.asmx") | Out-Null
$settings.Range("B8").Style = "Hyperlink"

$settings.Hyperlinks.Add($settings.Range("B13"), "https://online.immi.gov.au/lusc/login") | Out-Null
$settings.Range("B13").Style = "Hyperlink"

$settings.Range("A9").Select()

# ---------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

$constants.Range("A1").Hyperlinks.Delete()

# Row 25 was Email_SenderName -> becomes NoEmail_Message.
$constants.Range("A25").Value = "NoEmail_Message"
$constants.Range("B25").Value = "No new email or attachment found in mailbox to process."
$constants.Range("C25").Value = "Error message if no email or attachment found in mailbox"

# Row 26 was Email_SenderAddress -> becomes Email_MessageBody (text trimmed).
$constants.Range("A26").Value = "Email_MessageBody"
$constants.Range("B26").Value = "Please find VEVO location check result attached. There are {0} records processed. "
$constants.Range("B26").Style = "Normal"
$constants.Range("C26").Value = "Email message body for response email"

# Row 27 was Email_MessageBody -> becomes Email_AlertRecipient (hyperlink kept).
$constants.Range("A27").Value = "Email_AlertRecipient"
$constants.Range("B27").Value = "rpa.ads@its.uq.edu.au"
$constants.Range("C27").Value = "RPA ADS support email address   "

# Row 28 was Email_AlertRecipient -> becomes Email_Alert_Subject (now populated).
$constants.Range("A28").Value = "Email_Alert_Subject"
$constants.Range("B28").Value = "PBI_LocationCheck - Error Notification"
$constants.Range("B28").Style = "Hyperlink"
$constants.Range("C28").Value = "Email subject for error notification. "

# Row 29 was Email_Alert_Subject (empty) -> becomes EmailMsg_InitError.
$constants.Range("A29").Value = "EmailMsg_InitError"
$constants.Range("B29").Value = "Attention - PBI_LocationCheck encounter error during Initialization with following error. "
$constants.Range("C29").Value = "Email message body for move message to archive folder error. "

# Row 30 was Email_Alert_MessageBody -> becomes EmailMsg_MoveFolder (no description).
$constants.Range("A30").Value = "EmailMsg_MoveFolder"
$constants.Range("B30").Value = "Attention - PBI_LocationCheck process the records and sent updated file to PBI team however encountered an error while moving the email to `u{2018}PBI_Archive`u{2019} folder which need manual intervention prior to next run. Kindly move email `u{2018}{0}`u{2019} to PBI_Archive folder in RPA00001 Mailbox. "
$constants.Range("C30").ClearContents()

# Only one hyperlink survives, now anchored to B27 (still pointing at the
# original RPA00001 mailbox mailto: target).
$constants.Hyperlinks.Add($constants.Range("B27"), "mailto:rpa00001@uq.edu.au") | Out-Null
$constants.Range("B27").Style = "Hyperlink"

$constants.Range("B35").Select()

# ---------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# UQ_SMTP_SERVER / UQ_SMTP_PORT assets are removed entirely.
$assets.Range("A5:C6").ClearContents()

$assets.Range("A22").Select()

# ---------------------------------------------------------------------
# Constants becomes the active sheet/tab.
# ---------------------------------------------------------------------
$constants.Activate()
